$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 416.66666
$ws.Range("J2").Value = 425
$ws.Range("L2").Value = 425
$ws.Range("N2").Value = -651
$ws.Range("H9").Value = 68.09090999999999
$ws.Range("I9").Value = 48.666668
$ws.Range("J9").Value = 91.40000000000001
$ws.Range("K9").Value = 48.666668
$ws.Range("L9").Value = 91.40000000000001
$ws.Range("M9").Value = 120.333332
$ws.Range("N9").Value = -429.4
$ws.Range("H42").Value = 588.1429000000001
$ws.Range("I42").Value = 30
$ws.Range("J42").Value = 811.4
$ws.Range("K42").Value = 90
$ws.Range("L42").Value = 2434.2
$ws.Range("M42").Value = 140
$ws.Range("N42").Value = -2894.2
$ws.Range("H74").Value = 4346
$ws.Range("I74").Value = 3562.25
$ws.Range("K74").Value = 3562.25
$ws.Range("M74").Value = -2626.25
$ws.Range("H77").Value = 4346
$ws.Range("I77").Value = 3562.25
$ws.Range("K77").Value = 17811.25
$ws.Range("M77").Value = -13131.25
$ws.Range("H98").Value = 1435.4131
$ws.Range("I98").Value = 1398.2683
$ws.Range("K98").Value = 1398.2683
$ws.Range("M98").Value = 99.73170000000005
$ws.Range("H106").Value = 3190.6875
$ws.Range("I106").Value = 3023.3076
$ws.Range("K106").Value = 3023.3076
$ws.Range("M106").Value = -2392.3076
$ws.Range("H113").Value = 74874.89999999999
$ws.Range("I113").Value = 7155.8
$ws.Range("J113").Value = 142594
$ws.Range("K113").Value = 7155.8
$ws.Range("L113").Value = 142594
$ws.Range("M113").Value = -3901.8
$ws.Range("N113").Value = -149102
$ws.Range("H116").Value = 15311302
$ws.Range("I116").Value = 20609532
$ws.Range("J116").Value = 5301.5557
$ws.Range("K116").Value = 20609532
$ws.Range("L116").Value = 5301.5557
$ws.Range("M116").Value = -20606090
$ws.Range("N116").Value = -12185.5557
$ws.Range("H122").Value = 1435.4131
$ws.Range("I122").Value = 1398.2683
$ws.Range("K122").Value = 4194.8049
$ws.Range("M122").Value = -1744.8049

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2294.36
$ws.Range("I32").Value = 2294.36
$ws.Range("K32").Value = 2294.36
$ws.Range("M32").Value = -2007.36
$ws.Range("H74").Value = 22728694
$ws.Range("I74").Value = 27778738
$ws.Range("J74").Value = 3499.5
$ws.Range("K74").Value = 27778738
$ws.Range("L74").Value = 3499.5
$ws.Range("M74").Value = -27777864
$ws.Range("N74").Value = -5247.5
$ws.Range("H77").Value = 22728694
$ws.Range("I77").Value = 27778738
$ws.Range("J77").Value = 3499.5
$ws.Range("K77").Value = 138893690
$ws.Range("L77").Value = 17497.5
$ws.Range("M77").Value = -138889322
$ws.Range("N77").Value = -26233.5
$ws.Range("H106").Value = 74981.664
$ws.Range("J106").Value = 74981.664
$ws.Range("L106").Value = 74981.664
$ws.Range("N106").Value = -77505.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 108874.21
$ws.Range("J132").Value = 108874.21
$ws.Range("L132").Value = 108874.21
$ws.Range("N132").Value = -118994.21
$ws.Range("H134").Value = 2665
$ws.Range("I134").Value = 2634.342
$ws.Range("K134").Value = 7903.026
$ws.Range("M134").Value = -5368.026

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 650
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H31").Value = 30306428
$ws.Range("I31").Value = 47620804
$ws.Range("J31").Value = 6271.3335
$ws.Range("K31").Value = 47620804
$ws.Range("L31").Value = 6271.3335
$ws.Range("M31").Value = -47620509
$ws.Range("N31").Value = -6861.3335
$ws.Range("H34").Value = 30306428
$ws.Range("I34").Value = 47620804
$ws.Range("J34").Value = 6271.3335
$ws.Range("K34").Value = 47620804
$ws.Range("L34").Value = 6271.3335
$ws.Range("M34").Value = -47620602
$ws.Range("N34").Value = -6675.3335
$ws.Range("H58").Value = 2136.8
$ws.Range("I58").Value = 1976.1428
$ws.Range("J58").Value = 2779.4285
$ws.Range("K58").Value = 1976.1428
$ws.Range("L58").Value = 2779.4285
$ws.Range("M58").Value = -1773.1428
$ws.Range("N58").Value = -3185.4285
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51498
$ws.Range("H69").Value = 14845.5
$ws.Range("I69").Value = 3717.4443
$ws.Range("K69").Value = 3717.4443
$ws.Range("M69").Value = -2968.4443
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -157488
$ws.Range("H72").Value = 14845.5
$ws.Range("I72").Value = 3717.4443
$ws.Range("K72").Value = 11152.3329
$ws.Range("M72").Value = -7408.332900000001
$ws.Range("H74").Value = 79999.664
$ws.Range("J74").Value = 79999.664
$ws.Range("L74").Value = 79999.664
$ws.Range("N74").Value = -81747.664
$ws.Range("H77").Value = 79999.664
$ws.Range("J77").Value = 79999.664
$ws.Range("L77").Value = 239998.992
$ws.Range("N77").Value = -248734.992
$ws.Range("H132").Value = 31748152
$ws.Range("I132").Value = 38097116
$ws.Range("J132").Value = 3341.4285
$ws.Range("K132").Value = 114291348
$ws.Range("L132").Value = 10024.2855
$ws.Range("M132").Value = -114288818
$ws.Range("N132").Value = -15084.2855
$ws.Range("H136").Value = 2136.8
$ws.Range("I136").Value = 1976.1428
$ws.Range("J136").Value = 2779.4285
$ws.Range("K136").Value = 5928.428400000001
$ws.Range("L136").Value = 8338.2855
$ws.Range("M136").Value = -3378.428400000001
$ws.Range("N136").Value = -13438.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 81360.46000000001
$ws.Range("I138").Value = 95420
$ws.Range("K138").Value = 286260
$ws.Range("M138").Value = -281120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 43989.5
$ws.Range("J53").Value = 43989.5
$ws.Range("L53").Value = 43989.5
$ws.Range("N53").Value = -45251.5
$ws.Range("H80").Value = 81622.234
$ws.Range("I80").Value = 129936.25
$ws.Range("J80").Value = 4319.8
$ws.Range("K80").Value = 129936.25
$ws.Range("L80").Value = 4319.8
$ws.Range("M80").Value = -128938.25
$ws.Range("N80").Value = -6315.8
$ws.Range("H83").Value = 81622.234
$ws.Range("I83").Value = 129936.25
$ws.Range("J83").Value = 4319.8
$ws.Range("K83").Value = 649681.25
$ws.Range("L83").Value = 21599
$ws.Range("M83").Value = -644689.25
$ws.Range("N83").Value = -31583
$ws.Range("H95").Value = 28624
$ws.Range("J95").Value = 28624
$ws.Range("L95").Value = 28624
$ws.Range("N95").Value = -34116
$ws.Range("H113").Value = 1485.6111
$ws.Range("I113").Value = 1454.6
$ws.Range("K113").Value = 1454.6
$ws.Range("M113").Value = 715.4000000000001
$ws.Range("H123").Value = 48019.332
$ws.Range("J123").Value = 48019.332
$ws.Range("L123").Value = 48019.332
$ws.Range("N123").Value = -52919.332
$ws.Range("H136").Value = 34695.258
$ws.Range("J136").Value = 34695.258
$ws.Range("L136").Value = 104085.774
$ws.Range("N136").Value = -109185.774
$ws.Range("H140").Value = 80032
$ws.Range("J140").Value = 94693.5
$ws.Range("L140").Value = 94693.5
$ws.Range("N140").Value = -105053.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 630.9524
$ws.Range("I55").Value = 243.875
$ws.Range("J55").Value = 869.1539
$ws.Range("K55").Value = 243.875
$ws.Range("L55").Value = 869.1539
$ws.Range("M55").Value = -70.875
$ws.Range("N55").Value = -1215.1539
$ws.Range("H68").Value = 7140.1333
$ws.Range("J68").Value = 9393.4
$ws.Range("L68").Value = 9393.4
$ws.Range("N68").Value = -10891.4
$ws.Range("H71").Value = 7140.1333
$ws.Range("J71").Value = 9393.4
$ws.Range("L71").Value = 46967
$ws.Range("N71").Value = -54455
$ws.Range("H136").Value = 3618.2
$ws.Range("I136").Value = 2572.75
$ws.Range("K136").Value = 7718.25
$ws.Range("M136").Value = -5168.25
$ws.Range("H140").Value = 67055.36
$ws.Range("J140").Value = 67055.36
$ws.Range("L140").Value = 67055.36
$ws.Range("N140").Value = -77415.36

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16035
$ws.Range("J81").Value = 17602.385
$ws.Range("L81").Value = 35204.77
$ws.Range("N81").Value = -37326.77
$ws.Range("H84").Value = 16035
$ws.Range("J84").Value = 17602.385
$ws.Range("L84").Value = 176023.85
$ws.Range("N84").Value = -186631.85
$ws.Range("H113").Value = 647.3684
$ws.Range("J113").Value = 1097.25
$ws.Range("L113").Value = 3291.75
$ws.Range("N113").Value = -7631.75
$ws.Range("H132").Value = 12349053
$ws.Range("I132").Value = 25643588
$ws.Range("J132").Value = 4128
$ws.Range("K132").Value = 76930764
$ws.Range("L132").Value = 12384
$ws.Range("M132").Value = -76928234
$ws.Range("N132").Value = -17444
$ws.Range("H136").Value = 3299.8433
$ws.Range("I136").Value = 2095.131
$ws.Range("J136").Value = 6640.1816
$ws.Range("K136").Value = 6285.393
$ws.Range("L136").Value = 19920.5448
$ws.Range("M136").Value = -3735.393
$ws.Range("N136").Value = -25020.5448
